$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.412.13"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.569.22"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.64"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3740"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.25"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3359"
$ws.Range("E9").Value = "  -2.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.127"
$ws.Range("E10").Value = "  -3.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07435"
$ws.Range("E11").Value = "  -3.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.91"
$ws.Range("E13").Value = "  -2.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.899"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.851"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").Value = "1.569.51"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001116"
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.93"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06676"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.156"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.16"
$ws.Range("E22").Value = "  -2.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.84"
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").Value = "22.399.44"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.367"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.527"
$ws.Range("E26").Value = "  -9.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.93"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "146.90"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.988"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.58"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").Value = "1.740.49"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.001"
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.966"
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.901"
$ws.Range("E34").Value = "  -5.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.674"
$ws.Range("E35").Value = "  -4.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08398"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.368"
$ws.Range("E37").Value = "  +3.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02442"
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2240"
$ws.Range("E39").Value = "  -3.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06387"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.377"
$ws.Range("E41").Value = "  -3.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.16"
$ws.Range("E42").Value = "  -5.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6198"
$ws.Range("E43").Value = "  -3.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.93"
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.802"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5787"
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.052"
$ws.Range("E48").Value = "  -2.28%  "
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.224"
$ws.Range("E50").Value = "  -3.84%  "
$ws.Range("E51").Value = "  +0.26%  "
